$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The dataset rows 2-5 (columns D and L:T) are cyclically shifted up by one row:
# new row2 <- old row3, new row3 <- old row4, new row4 <- old row5, new row5 <- old row2
# Capture the "before" values first so later writes don't clobber data we still need to read.

$rows = @(2, 3, 4, 5)
$cols = @(4, 12, 13, 14, 15, 16, 17, 18, 19, 20)  # D, L, M, N, O, P, Q, R, S, T

$snapshot = @{}
foreach ($r in $rows) {
    $rowData = @{}
    foreach ($c in $cols) {
        $rowData[$c] = $ws.Cells.Item($r, $c).Value()
    }
    $snapshot[$r] = $rowData
}

# Map: target row -> source row (cyclic shift up by one, row2 wraps from row5)
$mapping = @{ 2 = 3; 3 = 4; 4 = 5; 5 = 2 }

foreach ($targetRow in $rows) {
    $sourceRow = $mapping[$targetRow]
    $srcData = $snapshot[$sourceRow]
    foreach ($c in $cols) {
        $ws.Cells.Item($targetRow, $c).Value = $srcData[$c]
    }
}
